$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8201676528599606
$ws.Range("C2").Value = 0.1432348664414804

$ws.Range("B3").Value = 0.7920611439842209
$ws.Range("C3").Value = 0.1399943064866937

$ws.Range("B4").Value = 0.8178994082840237
$ws.Range("C4").Value = 0.1766325669289702

$ws.Range("B5").Value = 0.6530078895463511
$ws.Range("C5").Value = 0.1603199767724667

$ws.Range("B6").Value = 0.8366494082840237
$ws.Range("C6").Value = 0.08090648525212538
